$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- Row 62: Junkyard Planet ---
$ws.Range("A62").Value = "Junkyard Planet"
$ws.Range("B62").Value = "Adam Minter"
$ws.Range("C61").Copy($ws.Range("C62"))
$ws.Range("C62").Value = 43939
$ws.Range("D61").Copy($ws.Range("D62"))
$ws.Range("D62").Value = 43946
$ws.Range("E62").Value = "business;scrap;junk;recycling"
$ws.Range("F62").Value = "Ebook"
$ws.Range("G62").Value = "304 Pages"

# --- Row 63: The Organized Mind ---
$ws.Range("A63").Value = "The Organized Mind"
$ws.Range("B63").Value = "Daniel Levitin"
$ws.Range("C61").Copy($ws.Range("C63"))
$ws.Range("C63").Value = 43945
$ws.Range("D61").Copy($ws.Range("D63"))
$ws.Range("D63").Value = 43947
$ws.Range("E63").Value = "multi-tasking;memory;productivity;self-improvement;thinking;psychology"
$ws.Range("F63").Value = "Audio"
$ws.Range("G63").Value = "16 Hours 20 Mins"

# --- Update view: scroll & active cell to mirror the author's final position ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
[void]$ws.Range("A64").Select()

Write-Host "Added 'Junkyard Planet' and 'The Organized Mind' rows."
